$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): I0, IF ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match formatting of the existing header cells (bold, centered/top-aligned, boxed border)
$headerRng = $ws.Range("I1:J1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# --- New data values for columns I (I0) and J (IF), rows 2-27 ---
$values = @{
    2  = @(1, 5)
    3  = @(1, 6)
    4  = @(1, 6)
    5  = @(1, 4)
    6  = @(1, 7)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(1, 4)
    12 = @(1, 6)
    13 = @(1, 5)
    14 = @(1, 4)
    15 = @(1, 6)
    16 = @(1, 5)
    17 = @(1, 2)
    18 = @(1, 6)
    19 = @(1, 7)
    20 = @(1, 6)
    21 = @(1, 5)
    22 = @(1, 7)
    23 = @(1, 5)
    24 = @(1, 6)
    25 = @(1, 5)
    26 = @(6, 9)
    27 = @(4, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
